$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.865.47"
$ws.Range("E2").Value = "  -0.40%  "

$ws.Range("D3").Value = "2.364.22"
$ws.Range("E3").Value = "  +0.18%  "

$ws.Range("E4").Value = "  -0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.672"
$ws.Range("E5").Value = "  -1.44%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "240.06"
$ws.Range("E6").Value = "  +0.10%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "74.56"
$ws.Range("E7").Value = "  +0.67%  "

$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("E9").Value = "  +0.90%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "60.20"
$ws.Range("E11").Value = "  +5.10%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "37.15"
$ws.Range("E12").Value = "  +15.15%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.28"
$ws.Range("E13").Value = "  +0.11%  "

$ws.Range("B14").Value = "TRON"
$ws.Range("C14").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.108"
$ws.Range("E14").Value = "  +0.46%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "16.41"
$ws.Range("E15").Value = "  -1.00%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.933"
$ws.Range("E16").Value = "  +3.22%  "

$ws.Range("D17").Value = "2.368.99"
$ws.Range("E17").Value = "  +0.22%  "

$ws.Range("D18").Value = "43.833.24"
$ws.Range("E18").Value = "  -0.31%  "

$ws.Range("E19").Value = "  +2.38%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.62"
$ws.Range("E20").Value = "  -4.82%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "77.36"
$ws.Range("E21").Value = "  +0.08%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "251.54"
$ws.Range("E22").Value = "  -2.91%  "

$ws.Range("E23").Value = "  +3.64%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("E24").Value = "  -0.06%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.87"
$ws.Range("E25").Value = "  -6.09%  "

$ws.Range("E26").Value = "  +0.14%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.55"
$ws.Range("E27").Value = "  -2.41%  "

$ws.Range("E28").Value = "  +2.73%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "22.41"
$ws.Range("E29").Value = "  -1.70%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "175.31"
$ws.Range("E30").Value = "  -0.15%  "

$ws.Range("E31").Value = "  +0.89%  "

$ws.Range("E32").Value = "  -2.04%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0759"
$ws.Range("E33").Value = "  -0.34%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.47"
$ws.Range("E34").Value = "  -1.98%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.11"
$ws.Range("E35").Value = "  -2.41%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.80"
$ws.Range("E36").Value = "  +1.12%  "

$ws.Range("E37").Value = "  +3.86%  "

$ws.Range("E38").Value = "  +1.89%  "

$ws.Range("E39").Value = "  +0.31%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.62"
$ws.Range("E40").Value = "  +18.31%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "20.85"
$ws.Range("E41").Value = "  +9.63%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "65.12"
$ws.Range("E42").Value = "  +10.92%  "

$ws.Range("E43").Value = "  -4.52%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "9.09"
$ws.Range("E44").Value = "  +0.96%  "

$ws.Range("E45").Value = "  -0.41%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.54"
$ws.Range("E46").Value = "  +1.60%  "

$ws.Range("E47").Value = "  -0.05%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.999"
$ws.Range("E48").Value = "  -0.26%  "

$ws.Range("E49").Value = "  -1.45%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "98.59"
$ws.Range("E50").Value = "  -2.27%  "

$ws.Range("E51").Value = "  +2.26%  "

